# Models Used.xlsx - add "Possum Morphometric Analysis" project rows
# (Random Forest Regressor / Non-Linear Regression) and
# (Random Forest Classifier / Binary Classification), inserted right
# after the existing "Random Forest Regressor / Laptop Price Analysis" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new blank rows at 7:8 (everything from old row 7 downward
#    shifts down to row 9+).
$ws.Rows("7:8").Insert()

# 2. Copy the formatting of row 6 (Random Forest Regressor / Non-Linear
#    Regression / Laptop Price Analysis) onto the two new rows - this is
#    the same visual style the new rows should use.
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C8").PasteSpecial(-4122)

# 3. Fill in the new cell values. Entering "Random Forest Classifier" then
#    "Possum Morphometric Analysis" first (before the other, already-used
#    strings) matches the shared-string insertion order of the target file.
$ws.Range("A8").Value = "Random Forest Classifier"
$ws.Range("C7").Value = "Possum Morphometric Analysis"
$ws.Range("A7").Value = "Random Forest Regressor"
$ws.Range("B7").Value = "Non-Linear Regression"
$ws.Range("B8").Value = "Binary Classification"
$ws.Range("C8").Value = "Possum Morphometric Analysis"

# 4. The engine does not auto-shift Hyperlink ranges when rows are
#    inserted, so rebuild the whole Hyperlinks collection pointing at the
#    correct (now shifted) cells, in the exact order that reproduces the
#    relationship-id numbering of the target workbook.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C3"), "c. Jupyter Notebooks\Bike Rental Demand.ipynb")
$ws.Hyperlinks.Add($ws.Range("C4"), "c. Jupyter Notebooks\Temperature Trends.ipynb")
$ws.Hyperlinks.Add($ws.Range("C9"), "c. Jupyter Notebooks\Graduate Admissions.ipynb")
$ws.Hyperlinks.Add($ws.Range("C11"), "c. Jupyter Notebooks\Employee Attrition.ipynb")
$ws.Hyperlinks.Add($ws.Range("C12"), "c. Jupyter Notebooks\Cardiac Risk.ipynb")
$ws.Hyperlinks.Add($ws.Range("C13"), "c. Jupyter Notebooks\Bank Churn.ipynb")
$ws.Hyperlinks.Add($ws.Range("C14"), "c. Jupyter Notebooks\Skin Analysis.ipynb")
$ws.Hyperlinks.Add($ws.Range("C16"), "c. Jupyter Notebooks\Air Passenger.ipynb")
$ws.Hyperlinks.Add($ws.Range("C17"), "c. Jupyter Notebooks\Champagne Sales.ipynb")
$ws.Hyperlinks.Add($ws.Range("C5"), "c. Jupyter Notebooks\King County House Sales.ipynb")
$ws.Hyperlinks.Add($ws.Range("C10"), "c. Jupyter Notebooks\Bank Client Term Deposit.ipynb")
$ws.Hyperlinks.Add($ws.Range("C15"), "c. Jupyter Notebooks\National Health and Nutrition Examination Survey (NHANES).ipynb", "", "", "NHANES")
$ws.Range("C15").Value = "National Health and Nutrition Examination Survey (NHANES)"
$ws.Hyperlinks.Add($ws.Range("C19"), "c. Jupyter Notebooks\Oil Prices.ipynb")
$ws.Hyperlinks.Add($ws.Range("C6"), "c. Jupyter Notebooks\Laptop Price Analysis.ipynb")
$ws.Hyperlinks.Add($ws.Range("C18"), "c. Jupyter Notebooks\Mindtree Stock Price.ipynb")
$ws.Hyperlinks.Add($ws.Range("C7"), "c. Jupyter Notebooks\Possum Morphometric Analysis.ipynb")
$ws.Hyperlinks.Add($ws.Range("C8"), "c. Jupyter Notebooks\Possum Morphometric Analysis.ipynb")

# 5. Update the saved view: scroll back to the top and select A5 (matches
#    the committed file - the "top-left cell" pin is cleared automatically
#    by selecting a cell that is already on screen).
$ws.Range("A5").Select()
